$wb = $excel.ActiveWorkbook

# The workbook has two sheets that mirror the same event data: "展览" and "全部类型".
# Update the "想去人数" (want-to-go count) column (F) for the rows that changed.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3223
    $ws.Range("F4").Value = 1099
}
